$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Itgav"
$ws.Cells.Item(2, 3).Value = "Thy1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 18.382477
$ws.Cells.Item(2, 8).Value = 55.147431
$ws.Cells.Item(2, 9).Value = 0.06380158579420245
$ws.Cells.Item(2, 10).Value = 0.06380158579420243
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.389838666666666
$ws.Cells.Item(2, 14).Value = 10.169516
$ws.Cells.Item(2, 15).Value = 0.03393930434450846
$ws.Cells.Item(2, 16).Value = 0.03393930434450846
$ws.Cells.Item(2, 17).Value = 62.31363132371066
$ws.Cells.Item(2, 18).Value = 560.822681913396
$ws.Cells.Item(2, 19).Value = 0.002165381437931704
$ws.Cells.Item(2, 20).Value = 0.002165381437931704

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Itgav"
$ws.Cells.Item(3, 3).Value = "Thy1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 18.382477
$ws.Cells.Item(3, 8).Value = 55.147431
$ws.Cells.Item(3, 9).Value = 0.06380158579420245
$ws.Cells.Item(3, 10).Value = 0.06380158579420243
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 84.50377666666667
$ws.Cells.Item(3, 14).Value = 253.51133
$ws.Cells.Item(3, 15).Value = 0.8460577852132902
$ws.Cells.Item(3, 16).Value = 0.8460577852132902
$ws.Cells.Item(3, 17).Value = 1553.388730988137
$ws.Cells.Item(3, 18).Value = 13980.49857889323
$ws.Cells.Item(3, 19).Value = 0.05397982837013864
$ws.Cells.Item(3, 20).Value = 0.05397982837013863

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Itgav"
$ws.Cells.Item(4, 3).Value = "Thy1"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 18.382477
$ws.Cells.Item(4, 8).Value = 55.147431
$ws.Cells.Item(4, 9).Value = 0.06380158579420245
$ws.Cells.Item(4, 10).Value = 0.06380158579420243
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1070423333333333
$ws.Cells.Item(4, 14).Value = 0.321127
$ws.Cells.Item(4, 15).Value = 0.001071715407718417
$ws.Cells.Item(4, 16).Value = 0.001071715407718417
$ws.Cells.Item(4, 17).Value = 1.967703230526333
$ws.Cells.Item(4, 18).Value = 17.709329074737
$ws.Cells.Item(4, 19).Value = 0.00006837714253251526
$ws.Cells.Item(4, 20).Value = 0.00006837714253251525

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Itgav"
$ws.Cells.Item(5, 3).Value = "Thy1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 18.382477
$ws.Cells.Item(5, 8).Value = 55.147431
$ws.Cells.Item(5, 9).Value = 0.06380158579420245
$ws.Cells.Item(5, 10).Value = 0.06380158579420243
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 11.878781
$ws.Cells.Item(5, 14).Value = 35.636343
$ws.Cells.Item(5, 15).Value = 0.1189311950344829
$ws.Cells.Item(5, 16).Value = 0.1189311950344828
$ws.Cells.Item(5, 17).Value = 218.361418520537
$ws.Cells.Item(5, 18).Value = 1965.252766684833
$ws.Cells.Item(5, 19).Value = 0.007587998843599582
$ws.Cells.Item(5, 20).Value = 0.007587998843599579

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Itgav"
$ws.Cells.Item(6, 3).Value = "Thy1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 82.64333833333332
$ws.Cells.Item(6, 8).Value = 247.930015
$ws.Cells.Item(6, 9).Value = 0.2868370808239535
$ws.Cells.Item(6, 10).Value = 0.2868370808239535
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.389838666666666
$ws.Cells.Item(6, 14).Value = 10.169516
$ws.Cells.Item(6, 15).Value = 0.03393930434450846
$ws.Cells.Item(6, 16).Value = 0.03393930434450846
$ws.Cells.Item(6, 17).Value = 280.1475838247488
$ws.Cells.Item(6, 18).Value = 2521.32825442274
$ws.Cells.Item(6, 19).Value = 0.009735050983374528
$ws.Cells.Item(6, 20).Value = 0.009735050983374528

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Itgav"
$ws.Cells.Item(7, 3).Value = "Thy1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 82.64333833333332
$ws.Cells.Item(7, 8).Value = 247.930015
$ws.Cells.Item(7, 9).Value = 0.2868370808239535
$ws.Cells.Item(7, 10).Value = 0.2868370808239535
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 84.50377666666667
$ws.Cells.Item(7, 14).Value = 253.51133
$ws.Cells.Item(7, 15).Value = 0.8460577852132902
$ws.Cells.Item(7, 16).Value = 0.8460577852132902
$ws.Cells.Item(7, 17).Value = 6983.674205507771
$ws.Cells.Item(7, 18).Value = 62853.06784956995
$ws.Cells.Item(7, 19).Value = 0.2426807453189596
$ws.Cells.Item(7, 20).Value = 0.2426807453189596

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Itgav"
$ws.Cells.Item(8, 3).Value = "Thy1"
$ws.Cells.Item(8, 4).Value = "M1"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 82.64333833333332
$ws.Cells.Item(8, 8).Value = 247.930015
$ws.Cells.Item(8, 9).Value = 0.2868370808239535
$ws.Cells.Item(8, 10).Value = 0.2868370808239535
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1070423333333333
$ws.Cells.Item(8, 14).Value = 0.321127
$ws.Cells.Item(8, 15).Value = 0.001071715407718417
$ws.Cells.Item(8, 16).Value = 0.001071715407718417
$ws.Cells.Item(8, 17).Value = 8.846335769656109
$ws.Cells.Item(8, 18).Value = 79.61702192690498
$ws.Cells.Item(8, 19).Value = 0.000307407719024004
$ws.Cells.Item(8, 20).Value = 0.000307407719024004

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Itgav"
$ws.Cells.Item(9, 3).Value = "Thy1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 82.64333833333332
$ws.Cells.Item(9, 8).Value = 247.930015
$ws.Cells.Item(9, 9).Value = 0.2868370808239535
$ws.Cells.Item(9, 10).Value = 0.2868370808239535
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 11.878781
$ws.Cells.Item(9, 14).Value = 35.636343
$ws.Cells.Item(9, 15).Value = 0.1189311950344829
$ws.Cells.Item(9, 16).Value = 0.1189311950344828
$ws.Cells.Item(9, 17).Value = 981.7021171705717
$ws.Cells.Item(9, 18).Value = 8835.319054535144
$ws.Cells.Item(9, 19).Value = 0.03411387680259533
$ws.Cells.Item(9, 20).Value = 0.03411387680259533

$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Itgav"
$ws.Cells.Item(10, 3).Value = "Thy1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 91.25099666666667
$ws.Cells.Item(10, 8).Value = 273.75299
$ws.Cells.Item(10, 9).Value = 0.3167123936907314
$ws.Cells.Item(10, 10).Value = 0.3167123936907314
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.389838666666666
$ws.Cells.Item(10, 14).Value = 10.169516
$ws.Cells.Item(10, 15).Value = 0.03393930434450846
$ws.Cells.Item(10, 16).Value = 0.03393930434450846
$ws.Cells.Item(10, 17).Value = 309.3261568725378
$ws.Cells.Item(10, 18).Value = 2783.93541185284
$ws.Cells.Item(10, 19).Value = 0.01074899831914751
$ws.Cells.Item(10, 20).Value = 0.01074899831914751

$ws.Cells.Item(11, 1).Value = "M1"
$ws.Cells.Item(11, 2).Value = "Itgav"
$ws.Cells.Item(11, 3).Value = "Thy1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 91.25099666666667
$ws.Cells.Item(11, 8).Value = 273.75299
$ws.Cells.Item(11, 9).Value = 0.3167123936907314
$ws.Cells.Item(11, 10).Value = 0.3167123936907314
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 84.50377666666667
$ws.Cells.Item(11, 14).Value = 253.51133
$ws.Cells.Item(11, 15).Value = 0.8460577852132902
$ws.Cells.Item(11, 16).Value = 0.8460577852132902
$ws.Cells.Item(11, 17).Value = 7711.053842930744
$ws.Cells.Item(11, 18).Value = 69399.48458637671
$ws.Cells.Item(11, 19).Value = 0.2679569863555799
$ws.Cells.Item(11, 20).Value = 0.2679569863555799

$ws.Cells.Item(12, 1).Value = "M1"
$ws.Cells.Item(12, 2).Value = "Itgav"
$ws.Cells.Item(12, 3).Value = "Thy1"
$ws.Cells.Item(12, 4).Value = "M1"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 91.25099666666667
$ws.Cells.Item(12, 8).Value = 273.75299
$ws.Cells.Item(12, 9).Value = 0.3167123936907314
$ws.Cells.Item(12, 10).Value = 0.3167123936907314
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.1070423333333333
$ws.Cells.Item(12, 14).Value = 0.321127
$ws.Cells.Item(12, 15).Value = 0.001071715407718417
$ws.Cells.Item(12, 16).Value = 0.001071715407718417
$ws.Cells.Item(12, 17).Value = 9.767719602192223
$ws.Cells.Item(12, 18).Value = 87.90947641973
$ws.Cells.Item(12, 19).Value = 0.0003394255521337381
$ws.Cells.Item(12, 20).Value = 0.0003394255521337381

$ws.Cells.Item(13, 1).Value = "M1"
$ws.Cells.Item(13, 2).Value = "Itgav"
$ws.Cells.Item(13, 3).Value = "Thy1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 91.25099666666667
$ws.Cells.Item(13, 8).Value = 273.75299
$ws.Cells.Item(13, 9).Value = 0.3167123936907314
$ws.Cells.Item(13, 10).Value = 0.3167123936907314
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 11.878781
$ws.Cells.Item(13, 14).Value = 35.636343
$ws.Cells.Item(13, 15).Value = 0.1189311950344829
$ws.Cells.Item(13, 16).Value = 0.1189311950344828
$ws.Cells.Item(13, 17).Value = 1083.950605435063
$ws.Cells.Item(13, 18).Value = 9755.555448915571
$ws.Cells.Item(13, 19).Value = 0.03766698346387029
$ws.Cells.Item(13, 20).Value = 0.03766698346387029

$ws.Cells.Item(14, 1).Value = "M2"
$ws.Cells.Item(14, 2).Value = "Itgav"
$ws.Cells.Item(14, 3).Value = "Thy1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 81.28845566666666
$ws.Cells.Item(14, 8).Value = 243.865367
$ws.Cells.Item(14, 9).Value = 0.2821345773094157
$ws.Cells.Item(14, 10).Value = 0.2821345773094157
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.389838666666666
$ws.Cells.Item(14, 14).Value = 10.169516
$ws.Cells.Item(14, 15).Value = 0.03393930434450846
$ws.Cells.Item(14, 16).Value = 0.03393930434450846
$ws.Cells.Item(14, 17).Value = 275.5547501724857
$ws.Cells.Item(14, 18).Value = 2479.992751552372
$ws.Cells.Item(14, 19).Value = 0.009575451285413509
$ws.Cells.Item(14, 20).Value = 0.009575451285413509

$ws.Cells.Item(15, 1).Value = "M2"
$ws.Cells.Item(15, 2).Value = "Itgav"
$ws.Cells.Item(15, 3).Value = "Thy1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 81.28845566666666
$ws.Cells.Item(15, 8).Value = 243.865367
$ws.Cells.Item(15, 9).Value = 0.2821345773094157
$ws.Cells.Item(15, 10).Value = 0.2821345773094157
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 84.50377666666667
$ws.Cells.Item(15, 14).Value = 253.51133
$ws.Cells.Item(15, 15).Value = 0.8460577852132902
$ws.Cells.Item(15, 16).Value = 0.8460577852132902
$ws.Cells.Item(15, 17).Value = 6869.181503234235
$ws.Cells.Item(15, 18).Value = 61822.63352910811
$ws.Cells.Item(15, 19).Value = 0.238702155610492
$ws.Cells.Item(15, 20).Value = 0.238702155610492

$ws.Cells.Item(16, 1).Value = "M2"
$ws.Cells.Item(16, 2).Value = "Itgav"
$ws.Cells.Item(16, 3).Value = "Thy1"
$ws.Cells.Item(16, 4).Value = "M1"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 81.28845566666666
$ws.Cells.Item(16, 8).Value = 243.865367
$ws.Cells.Item(16, 9).Value = 0.2821345773094157
$ws.Cells.Item(16, 10).Value = 0.2821345773094157
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1070423333333333
$ws.Cells.Item(16, 14).Value = 0.321127
$ws.Cells.Item(16, 15).Value = 0.001071715407718417
$ws.Cells.Item(16, 16).Value = 0.001071715407718417
$ws.Cells.Item(16, 17).Value = 8.701305967623222
$ws.Cells.Item(16, 18).Value = 78.31175370860899
$ws.Cells.Item(16, 19).Value = 0.0003023679735526238
$ws.Cells.Item(16, 20).Value = 0.0003023679735526238

$ws.Cells.Item(17, 1).Value = "M2"
$ws.Cells.Item(17, 2).Value = "Itgav"
$ws.Cells.Item(17, 3).Value = "Thy1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 81.28845566666666
$ws.Cells.Item(17, 8).Value = 243.865367
$ws.Cells.Item(17, 9).Value = 0.2821345773094157
$ws.Cells.Item(17, 10).Value = 0.2821345773094157
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 11.878781
$ws.Cells.Item(17, 14).Value = 35.636343
$ws.Cells.Item(17, 15).Value = 0.1189311950344829
$ws.Cells.Item(17, 16).Value = 0.1189311950344828
$ws.Cells.Item(17, 17).Value = 965.6077626925425
$ws.Cells.Item(17, 18).Value = 8690.469864232882
$ws.Cells.Item(17, 19).Value = 0.03355460243995749
$ws.Cells.Item(17, 20).Value = 0.03355460243995749

$ws.Cells.Item(18, 1).Value = "sCs"
$ws.Cells.Item(18, 2).Value = "Itgav"
$ws.Cells.Item(18, 3).Value = "Thy1"
$ws.Cells.Item(18, 4).Value = "ECs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 14.55416966666667
$ws.Cells.Item(18, 8).Value = 43.662509
$ws.Cells.Item(18, 9).Value = 0.0505143623816971
$ws.Cells.Item(18, 10).Value = 0.0505143623816971
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 3.389838666666666
$ws.Cells.Item(18, 14).Value = 10.169516
$ws.Cells.Item(18, 15).Value = 0.03393930434450846
$ws.Cells.Item(18, 16).Value = 0.03393930434450846
$ws.Cells.Item(18, 17).Value = 49.33628709729378
$ws.Cells.Item(18, 18).Value = 444.026583875644
$ws.Cells.Item(18, 19).Value = 0.001714422318641207
$ws.Cells.Item(18, 20).Value = 0.001714422318641207

$ws.Cells.Item(19, 1).Value = "sCs"
$ws.Cells.Item(19, 2).Value = "Itgav"
$ws.Cells.Item(19, 3).Value = "Thy1"
$ws.Cells.Item(19, 4).Value = "FAPs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 14.55416966666667
$ws.Cells.Item(19, 8).Value = 43.662509
$ws.Cells.Item(19, 9).Value = 0.0505143623816971
$ws.Cells.Item(19, 10).Value = 0.0505143623816971
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 84.50377666666667
$ws.Cells.Item(19, 14).Value = 253.51133
$ws.Cells.Item(19, 15).Value = 0.8460577852132902
$ws.Cells.Item(19, 16).Value = 0.8460577852132902
$ws.Cells.Item(19, 17).Value = 1229.882303080774
$ws.Cells.Item(19, 18).Value = 11068.94072772697
$ws.Cells.Item(19, 19).Value = 0.0427380695581202
$ws.Cells.Item(19, 20).Value = 0.04273806955812019

$ws.Cells.Item(20, 1).Value = "sCs"
$ws.Cells.Item(20, 2).Value = "Itgav"
$ws.Cells.Item(20, 3).Value = "Thy1"
$ws.Cells.Item(20, 4).Value = "M1"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 14.55416966666667
$ws.Cells.Item(20, 8).Value = 43.662509
$ws.Cells.Item(20, 9).Value = 0.0505143623816971
$ws.Cells.Item(20, 10).Value = 0.0505143623816971
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 0.6666666666666666
$ws.Cells.Item(20, 13).Value = 0.1070423333333333
$ws.Cells.Item(20, 14).Value = 0.321127
$ws.Cells.Item(20, 15).Value = 0.001071715407718417
$ws.Cells.Item(20, 16).Value = 0.001071715407718417
$ws.Cells.Item(20, 17).Value = 1.557912280849222
$ws.Cells.Item(20, 18).Value = 14.021210527643
$ws.Cells.Item(20, 19).Value = 0.0000541370204755364
$ws.Cells.Item(20, 20).Value = 0.00005413702047553639

$ws.Cells.Item(21, 1).Value = "sCs"
$ws.Cells.Item(21, 2).Value = "Itgav"
$ws.Cells.Item(21, 3).Value = "Thy1"
$ws.Cells.Item(21, 4).Value = "sCs"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 14.55416966666667
$ws.Cells.Item(21, 8).Value = 43.662509
$ws.Cells.Item(21, 9).Value = 0.0505143623816971
$ws.Cells.Item(21, 10).Value = 0.0505143623816971
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 11.878781
$ws.Cells.Item(21, 14).Value = 35.636343
$ws.Cells.Item(21, 15).Value = 0.1189311950344829
$ws.Cells.Item(21, 16).Value = 0.1189311950344828
$ws.Cells.Item(21, 17).Value = 172.8857941071764
$ws.Cells.Item(21, 18).Value = 1555.972146964587
$ws.Cells.Item(21, 19).Value = 0.006007733484460162
$ws.Cells.Item(21, 20).Value = 0.006007733484460161

